$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 16613.154
$ws.Range("I18").Value = 15748.375
$ws.Range("K18").Value = 15748.375
$ws.Range("M18").Value = -15464.375
# Row 112
$ws.Range("H112").Value = 1907.44
$ws.Range("J112").Value = 1907.44
$ws.Range("L112").Value = 5722.32
$ws.Range("N112").Value = -7938.32
# Row 113
$ws.Range("H113").Value = 26875.375
$ws.Range("I113").Value = 41701
$ws.Range("J113").Value = 2166
$ws.Range("K113").Value = 41701
$ws.Range("L113").Value = 2166
$ws.Range("M113").Value = -38447
$ws.Range("N113").Value = -8674
# Row 137
$ws.Range("H137").Value = 46668.137
$ws.Range("I137").Value = 869.0769
$ws.Range("J137").Value = 112822.336
$ws.Range("K137").Value = 2607.2307
$ws.Range("L137").Value = 338467.008
$ws.Range("M137").Value = -57.23070000000007
$ws.Range("N137").Value = -343567.008

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 327249.28
$ws.Range("I2").Value = 370775
$ws.Range("J2").Value = 806.5
$ws.Range("K2").Value = 370775
$ws.Range("L2").Value = 806.5
$ws.Range("M2").Value = -370662
$ws.Range("N2").Value = -1032.5
# Row 32
$ws.Range("H32").Value = 3693.36
$ws.Range("I32").Value = 2345.318
$ws.Range("J32").Value = 13579
$ws.Range("K32").Value = 2345.318
$ws.Range("L32").Value = 13579
$ws.Range("M32").Value = -2058.318
$ws.Range("N32").Value = -14153
# Row 45
$ws.Range("H45").Value = 1277.8948
$ws.Range("I45").Value = 889.8182
$ws.Range("K45").Value = 889.8182
$ws.Range("M45").Value = -512.8182
# Row 61
$ws.Range("H61").Value = 28488.033
$ws.Range("I61").Value = 35885.477
$ws.Range("K61").Value = 35885.477
$ws.Range("M61").Value = -35673.477
# Row 102
$ws.Range("H102").Value = 1722.375
$ws.Range("I102").Value = 1722.375
$ws.Range("K102").Value = 1722.375
$ws.Range("M102").Value = -100.375
# Row 116
$ws.Range("H116").Value = 327249.28
$ws.Range("I116").Value = 370775
$ws.Range("J116").Value = 806.5
$ws.Range("K116").Value = 370775
$ws.Range("L116").Value = 806.5
$ws.Range("M116").Value = -368481
$ws.Range("N116").Value = -5394.5
# Row 132
$ws.Range("H132").Value = 1507.551
$ws.Range("I132").Value = 1012.78125
$ws.Range("J132").Value = 2438.8823
$ws.Range("K132").Value = 3038.34375
$ws.Range("L132").Value = 7316.646900000001
$ws.Range("M132").Value = -508.34375
$ws.Range("N132").Value = -12376.6469
# Row 136
$ws.Range("H136").Value = 28488.033
$ws.Range("I136").Value = 35885.477
$ws.Range("K136").Value = 107656.431
$ws.Range("M136").Value = -105106.431

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 327249.28
$ws.Range("I3").Value = 370775
$ws.Range("J3").Value = 806.5
$ws.Range("K3").Value = 370775
$ws.Range("L3").Value = 806.5
$ws.Range("M3").Value = -370661
$ws.Range("N3").Value = -1034.5
# Row 5
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 5000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 5000
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -5226
# Row 94
$ws.Range("H94").Value = 725.36
$ws.Range("I94").Value = 336.14285
$ws.Range("J94").Value = 2768.75
$ws.Range("K94").Value = 336.14285
$ws.Range("L94").Value = 2768.75
$ws.Range("M94").Value = 114.85715
$ws.Range("N94").Value = -3670.75
# Row 107
$ws.Range("H107").Value = 762.2
$ws.Range("I107").Value = 601.5
$ws.Range("J107").Value = 1003.25
$ws.Range("K107").Value = 601.5
$ws.Range("L107").Value = 1003.25
$ws.Range("M107").Value = 1318.5
$ws.Range("N107").Value = -4843.25
# Row 134
$ws.Range("H134").Value = 7934.2607
$ws.Range("I134").Value = 8823
$ws.Range("K134").Value = 26469
$ws.Range("M134").Value = -23934

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2068.0386
$ws.Range("I31").Value = 1638.95
$ws.Range("J31").Value = 3498.3333
$ws.Range("K31").Value = 1638.95
$ws.Range("L31").Value = 3498.3333
$ws.Range("M31").Value = -1343.95
$ws.Range("N31").Value = -4088.3333
# Row 34
$ws.Range("H34").Value = 2068.0386
$ws.Range("I34").Value = 1638.95
$ws.Range("J34").Value = 3498.3333
$ws.Range("K34").Value = 1638.95
$ws.Range("L34").Value = 3498.3333
$ws.Range("M34").Value = -1436.95
$ws.Range("N34").Value = -3902.3333
# Row 105
$ws.Range("H105").Value = 1191.6364
$ws.Range("I105").Value = 1200.8
$ws.Range("K105").Value = 1200.8
$ws.Range("M105").Value = 546.2
# Row 132
$ws.Range("H132").Value = 1571.7142
$ws.Range("I132").Value = 1276.9231
$ws.Range("J132").Value = 2050.75
$ws.Range("K132").Value = 3830.7693
$ws.Range("L132").Value = 6152.25
$ws.Range("M132").Value = -1300.7693
$ws.Range("N132").Value = -11212.25
# Row 134
$ws.Range("H134").Value = 2292.261
$ws.Range("I134").Value = 1974.8823
$ws.Range("J134").Value = 3191.5
$ws.Range("K134").Value = 5924.6469
$ws.Range("L134").Value = 9574.5
$ws.Range("M134").Value = -3389.6469
$ws.Range("N134").Value = -14644.5

$ws = $wb.Worksheets.Item("CUL")
# Row 44
$ws.Range("H44").Value = 2841.2
$ws.Range("I44").Value = 3768.6667
$ws.Range("J44").Value = 1450
$ws.Range("K44").Value = 11306.0001
$ws.Range("L44").Value = 4350
$ws.Range("M44").Value = -10908.0001
$ws.Range("N44").Value = -5146
# Row 131
$ws.Range("H131").Value = 17687.908
$ws.Range("J131").Value = 19070.34
$ws.Range("L131").Value = 57211.02
$ws.Range("N131").Value = -67291.02

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2499.8572
$ws.Range("I80").Value = 2124.75
$ws.Range("K80").Value = 2124.75
$ws.Range("M80").Value = -1126.75
# Row 83
$ws.Range("H83").Value = 2499.8572
$ws.Range("I83").Value = 2124.75
$ws.Range("K83").Value = 10623.75
$ws.Range("M83").Value = -5631.75
# Row 97
$ws.Range("H97").Value = 986.3
$ws.Range("I97").Value = 961.3333
$ws.Range("J97").Value = 1086.1666
$ws.Range("K97").Value = 961.3333
$ws.Range("L97").Value = 1086.1666
$ws.Range("M97").Value = -465.3333
$ws.Range("N97").Value = -2078.1666
# Row 102
$ws.Range("H102").Value = 3431.5
$ws.Range("I102").Value = 3035.875
$ws.Range("K102").Value = 3035.875
$ws.Range("M102").Value = -1413.875
# Row 113
$ws.Range("H113").Value = 1308.9333
$ws.Range("I113").Value = 973.3333
$ws.Range("J113").Value = 1812.3334
$ws.Range("K113").Value = 973.3333
$ws.Range("L113").Value = 1812.3334
$ws.Range("M113").Value = 1196.6667
$ws.Range("N113").Value = -6152.3334
# Row 139
$ws.Range("H139").Value = 70595
$ws.Range("J139").Value = 70595
$ws.Range("L139").Value = 70595
$ws.Range("N139").Value = -80875

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 5756.4165
$ws.Range("I16").Value = 6238.8184
$ws.Range("K16").Value = 6238.8184
$ws.Range("M16").Value = -6068.8184
# Row 93
$ws.Range("H93").Value = 967.44446
$ws.Range("I93").Value = 547.2174
$ws.Range("J93").Value = 3383.75
$ws.Range("K93").Value = 547.2174
$ws.Range("L93").Value = 3383.75
$ws.Range("M93").Value = 700.7826
$ws.Range("N93").Value = -5879.75
# Row 132
$ws.Range("H132").Value = 2114.5642
$ws.Range("I132").Value = 1527.9412
$ws.Range("K132").Value = 4583.8236
$ws.Range("M132").Value = -2053.8236

$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Range("H14").Value = 4500
$ws.Range("J14").Value = 4500
$ws.Range("L14").Value = 4500
$ws.Range("N14").Value = -4836
# Row 122
$ws.Range("H122").Value = 30061.678
$ws.Range("I122").Value = 39419.855
$ws.Range("K122").Value = 118259.565
$ws.Range("M122").Value = -115809.565
